# Rename the "Voros et al., 2005" references from mice-age labels to
# diet-duration labels on the two "adipose" sheets, then restore the
# selection / active-sheet state recorded in the saved workbook.

$wb = $excel.ActiveWorkbook

$wsVesselSizeAdipose    = $wb.Worksheets.Item("Vessel size (adipose)")
$wsVesselDensityAdipose = $wb.Worksheets.Item("Vessel density (adipose)")

# --- "Vessel size (adipose)" sheet ---------------------------------------
$wsVesselSizeAdipose.Range("A3").Value = "Voros et al., 2005 (2 wk. diet)"
$wsVesselSizeAdipose.Range("A4").Value = "Voros et al., 2005 (5 wk. diet)"
$wsVesselSizeAdipose.Range("A5").Value = "Voros et al., 2005 (15 wk. diet)"

# --- "Vessel density (adipose)" sheet ------------------------------------
$wsVesselDensityAdipose.Range("A3").Value = "Voros et al., 2005 (2 wk. diet)"
$wsVesselDensityAdipose.Range("A4").Value = "Voros et al., 2005 (5 wk. diet)"
$wsVesselDensityAdipose.Range("A5").Value = "Voros et al., 2005 (15 wk. diet)"

# --- Selection / active sheet bookkeeping --------------------------------
# "Vessel size (adipose)" is no longer the active tab; its selection moves
# to A3:A5 (still tracked even though it's not the foreground sheet).
[void]$wsVesselSizeAdipose.Range("A3:A5").Select()

# "Vessel density (adipose)" becomes the active/foreground tab with the
# same A3:A5 selection.
[void]$wsVesselDensityAdipose.Activate()
[void]$wsVesselDensityAdipose.Range("A3:A5").Select()
